# Insert a new price-record row for "Femacal de La Calera" (Poroto granado)
# right after the existing row 130, shifting all subsequent rows (131-165)
# down by one (they keep their original values). The newly inserted row 131
# gets its own date / price / origin data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 131..165 down to 132..166, leaving a blank row 131 to fill in.
$ws.Rows.Item(131).Insert()

$ws.Cells.Item(131, 1).Value  = 3
$ws.Cells.Item(131, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(131, 3).Value  = "Coquimbo"
$ws.Cells.Item(131, 4).Value  = 44642
$ws.Cells.Item(131, 5).Value  = 5
$ws.Cells.Item(131, 6).Value  = 100112030
$ws.Cells.Item(131, 7).Value  = "Poroto granado"
$ws.Cells.Item(131, 8).Value  = "Sin especificar"
$ws.Cells.Item(131, 9).Value  = "Primera"
$ws.Cells.Item(131, 10).Value = 73
$ws.Cells.Item(131, 11).Value = 21000
$ws.Cells.Item(131, 12).Value = 22000
$ws.Cells.Item(131, 13).Value = 21479
$ws.Cells.Item(131, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(131, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(131, 16).Value = 859
$ws.Cells.Item(131, 17).Value = 25
$ws.Cells.Item(131, 18).Value = "Hortaliza"
